# Road round endings and road bugfix
# Adds a new "Road circle" / mark row (row 29) to the 3D object library sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3d_object_library")

$row = 29

$ws.Cells.Item($row, 1).Value  = "23"
$ws.Cells.Item($row, 2).Value  = "Road circle"
$ws.Cells.Item($row, 3).Value  = "mark"
$ws.Cells.Item($row, 4).Value  = "road_circle.obj"
$ws.Cells.Item($row, 5).Value  = "*"
$ws.Cells.Item($row, 6).Value  = "road_circle.png"
$ws.Cells.Item($row, 7).Value  = "*"
$ws.Cells.Item($row, 8).Value  = "*"
$ws.Cells.Item($row, 9).Value  = "*"
$ws.Cells.Item($row, 10).Value = "*"

$ws.Cells.Item($row, 11).Value = "0.0"
$ws.Cells.Item($row, 12).Value = "0.0"
$ws.Cells.Item($row, 13).Value = "0.0"

$ws.Cells.Item($row, 14).Value = "1.0"
$ws.Cells.Item($row, 15).Value = "1.0"
$ws.Cells.Item($row, 16).Value = "1.0"

$ws.Cells.Item($row, 17).Value = "0.0"
$ws.Cells.Item($row, 18).Value = "0"
$ws.Cells.Item($row, 19).Value = "1"
$ws.Cells.Item($row, 20).Value = "0"

$ws.Cells.Item($row, 21).Value = "0.0"
$ws.Cells.Item($row, 22).Value = "0.0"
$ws.Cells.Item($row, 23).Value = "0.0"

$ws.Cells.Item($row, 24).Value = "2.0"

$ws.Cells.Item($row, 25).Value = "0.0"
$ws.Cells.Item($row, 26).Value = "0.0"
$ws.Cells.Item($row, 27).Value = "0.0"

$ws.Cells.Item($row, 28).Value = "2.0"

$ws.Cells.Item($row, 29).Value = "1.0"
$ws.Cells.Item($row, 30).Value = "1.0"
$ws.Cells.Item($row, 31).Value = "1.0"
$ws.Cells.Item($row, 32).Value = "1.0"

$ws.Cells.Item($row, 33).Value = "0"

# Move selection the way the authored workbook left it (below the newly
# added row), matching the saved view state in the target file.
$ws.Range("A30").Select()
